# Fixes total coal values in historical extension
#
# The Bond historical-extension "total coal" rows were mapped to the wrong
# ext_sector: Bond_to_ext!C23 and CEDS_to_ext!B2:B4 pointed at the existing
# " Power       " shared string (used elsewhere for genuinely power-sector
# rows) instead of a clean "Power" label. Re-enter the value as a fresh
# "Power" string so it gets its own shared-string entry distinct from the
# padded one.

$wb = $excel.ActiveWorkbook
$wsBond = $wb.Worksheets.Item("Bond_to_ext")
$wsCeds = $wb.Worksheets.Item("CEDS_to_ext")

# --- Data fix: re-point the coal-total rows at a clean "Power" string ----
$wsBond.Range("C23").Value = "Power"

$wsCeds.Range("B2").Value = "Power"
$wsCeds.Range("B3").Value = "Power"
$wsCeds.Range("B4").Value = "Power"

# --- View state: leave CEDS_to_ext as the active/selected tab, with the
#     last touched cells selected on each sheet -----------------------------
$wsBond.Activate()
$wsBond.Range("C23").Select()

$wsCeds.Activate()
$wsCeds.Range("A6").Select()
